$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.479.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "'2.108.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'333.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.5249"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("D8").Value = "'0.4571"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.91%  "
$ws.Range("D9").Value = "'53.46"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +13.76%  "
$ws.Range("D10").Value = "'0.08980"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").Value = "'24.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "'2.102.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "'6.806"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "'7.840"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'96.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").Value = "'0.00001131"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "'0.06618"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").Value = "'19.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").Value = "'6.305"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'30.551.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("D24").Value = "'12.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").Value = "'2.356"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("D26").Value = "'2.353.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'22.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("D28").Value = "'2.578"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").Value = "'163.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").Value = "'132.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").Value = "'1.705"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.19%  "
$ws.Range("D33").Value = "'0.1073"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'6.161"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "'3.932"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.34%  "
$ws.Range("E36").Value = "  +9.04%  "
$ws.Range("D37").Value = "'0.02577"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").Value = "'0.06828"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").Value = "'12.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("D41").Value = "'0.2291"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").Value = "'0.6908"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("D43").Value = "'1.243"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").Value = "'2.355"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.29%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "'14.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").Value = "'0.6388"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'3.651"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "'0.00000000356"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +25.17%  "
$ws.Range("D50").Value = "'1.248"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").Value = "'1.220"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.35%  "
